# v2.6 Added decoupled suspension, four-wheel steering, scripts to generate GGV diagram
#
# This script:
#  1) Adds a new "FSAE_Achilles" sheet (copied from Trailer_Kumanzi, the most
#     structurally-similar template) at the end of the workbook and fills in
#     its Aero-coefficient values.
#  2) Updates the CD value (H8) on Sedan_Hamba.
#  3) Refreshes the shared column widths on every sheet.
#  4) Restores view state (active sheet/selection) to match the authored
#     session.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) New sheet: FSAE_Achilles (copy Trailer_Kumanzi's layout/formatting)
# ---------------------------------------------------------------------------
$template = $wb.Worksheets.Item("Trailer_Kumanzi")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$template.Copy([System.Reflection.Missing]::Value, $lastSheet)

$fsae = $wb.Worksheets.Item($wb.Worksheets.Count)
$fsae.Name = "FSAE_Achilles"

# "class" label on the new sheet names itself
$fsae.Range("H3").Value = "FSAE_Achilles"

# Aero coefficients for the new vehicle
$fsae.Range("H5").Value = -2.5
$fsae.Range("H6").Value = 1
$fsae.Range("H8").Value = 1.2
$fsae.Range("F9").Value = -0.8
$fsae.Range("G9").Value = 0
$fsae.Range("H9").Value = 0.6

# ---------------------------------------------------------------------------
# 2) Sedan_Hamba: CD (H8) updated, now shown with a 2-decimal number format
# ---------------------------------------------------------------------------
$sedanHamba = $wb.Worksheets.Item("Sedan_Hamba")
$sedanHamba.Range("H8").Value = 1.98
$sedanHamba.Range("H8").NumberFormat = "0.00"

# ---------------------------------------------------------------------------
# 3) Shared column-width refresh across every sheet (incl. the new one)
# ---------------------------------------------------------------------------
$colWidths = @{ 1 = 13.666666666666666; 2 = 11.833333333333334; 3 = 14.833333333333334; 4 = 10.333333333333334; 5 = 15.833333333333334 }
foreach ($sheetName in @("Sedan_HambaLG","Sedan_Hamba","Bus_Makhulu","Truck_Amandla","Trailer_Elula","Trailer_Thwala","Trailer_Kumanzi","FSAE_Achilles")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Columns.Item(1).ColumnWidth = $colWidths[1]
    $ws.Columns.Item(2).ColumnWidth = $colWidths[2]
    $ws.Columns.Item(3).ColumnWidth = $colWidths[3]
    $ws.Columns.Item(4).ColumnWidth = $colWidths[4]
    $ws.Columns.Item(5).ColumnWidth = $colWidths[5]
    $ws.Range($ws.Columns.Item(9), $ws.Columns.Item(15)).ColumnWidth = 5.833333333333333
}

# ---------------------------------------------------------------------------
# 4) View state: active-cell selections that changed + final active sheet
# ---------------------------------------------------------------------------
$sedanHambaLG = $wb.Worksheets.Item("Sedan_HambaLG")
$sedanHambaLG.Activate()
$sedanHambaLG.Range("H8").Select()

$fsae.Activate()
$fsae.Range("G12").Select()

$sedanHamba.Activate()
$sedanHamba.Range("E18").Select()
